$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move existing F12 value ("_neck_clipped.obj") down to F13
$ws.Range("F13").Value = $ws.Range("F12").Value

# Update row 12 with the new patient-data prep values
$ws.Range("F12").Value = ".obj"
$ws.Range("G12").Value = "C:\Users\franz\Documents\work\projects\arp\data\patient_data\sagittal_patient_data_sept2023\sagittal_patient_data_sept2023_age_sex_data.xlsx"
$ws.Range("H12").Value = "C:\Users\franz\Documents\work\projects\arp\data\patient_data\sagittal_patient_data_sept2023\sagittal_patient_data_sept2023_age_sex_data.xlsx"
$ws.Range("I12").Value = $false
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = "_"
$ws.Range("M12").Value = $false
$ws.Range("N12").Value = $true
$ws.Range("O12").Value = $true

# Update view state to match the authored selection
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("L15").Select()
